$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-blank name cells
$ws.Range("B3").Value = "ARAUJO"
$ws.Range("B4").Value = "ARANDA"
$ws.Range("D5").Value = "BRENDA ALEJANDRA"
$ws.Range("D6").Value = "BRENDA PAOLA"
$ws.Range("D8").Value = "DAVID"

# Correct the "Integrantes familia" value on row 6
$ws.Range("E6").Value = 1

# Update the view: scroll over and change the active selection
$ws.Range("Z14").Select()
